$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9805734
Write-Host "done"
